# Prepend a new latest-quarter data point to the "Data" sheet.
# The sheet stores S&P 500 sales-growth figures with the most recent
# quarter in row 2 (row 1 is the header). A new row is inserted at row 2
# for 2025-09-30 (Excel serial 45930) with a value of 5.1, shifting all
# existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current first data row (row 2).
$ws.Rows.Item(2).Insert()

# Match the formatting of the other data rows (date format / centered
# number format) by copying the format from the row just below.
$ws.Range("A3:B3").Copy() | Out-Null
$ws.Range("A2:B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the new data point: 2025-09-30 -> 5.1% sales growth.
$ws.Range("A2").Value2 = 45930
$ws.Range("B2").Value2 = 5.1
